$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.00"
$ws.Range("E2").Value = "'-0.63%"
$ws.Range("E3").Value = "'-4.27%"
$ws.Range("D4").Value = "'5.247"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.05698"
$ws.Range("D6").Value = "'6.617"
$ws.Range("E6").Value = "'0.32%"
$ws.Range("D7").Value = "'3.191"
$ws.Range("E7").Value = "'3.39%"
$ws.Range("D8").Value = "'0.8504"
$ws.Range("E8").Value = "'-0.75%"
$ws.Range("D9").Value = "'0.8558"
$ws.Range("E9").Value = "'-1.49%"
$ws.Range("D10").Value = "'0.1368"
$ws.Range("E10").Value = "'0.52%"
$ws.Range("D11").Value = "'0.07068"
$ws.Range("E11").Value = "'-0.03%"
$ws.Range("D12").Value = "'0.03187"
$ws.Range("E12").Value = "'8.89%"
$ws.Range("E13").Value = "'-1.99%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.009987"
$ws.Range("E15").Value = "'1,556.80%"
$ws.Range("D16").Value = "'0.006037"
$ws.Range("E16").Value = "'-1.15%"
$ws.Range("D17").Value = "'3.491"
$ws.Range("E17").Value = "'0.11%"
$ws.Range("E18").Value = "'-4.09%"
$ws.Range("D19").Value = "'0.3169"
$ws.Range("E19").Value = "'0.42%"
$ws.Range("D20").Value = "'0.03268"
$ws.Range("E20").Value = "'-3.53%"
$ws.Range("D21").Value = "'0.1288"
$ws.Range("E21").Value = "'-2.20%"
$ws.Range("D22").Value = "'3.484"
$ws.Range("E22").Value = "'0.60%"
$ws.Range("D23").Value = "'0.04072"
$ws.Range("E23").Value = "'-2.42%"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("D26").Value = "'0.004139"
$ws.Range("E26").Value = "'-17.54%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.84%"
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("E40").Value = "'-0.09%"
$ws.Range("D41").Value = "'0.1062"
$ws.Range("E41").Value = "'-1.02%"
$ws.Range("D42").Value = "'0.003716"
$ws.Range("E42").Value = "'-35.23%"
$ws.Range("D43").Value = "'0.002400"
$ws.Range("E43").Value = "'20.05%"
$ws.Range("D44").Value = "'0.009336"
$ws.Range("E44").Value = "'-2.46%"
$ws.Range("D45").Value = "'0.00005275"
$ws.Range("E45").Value = "'1.01%"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("D47").Value = "'0.07496"
$ws.Range("E47").Value = "'15.89%"
$ws.Range("D48").Value = "'0.002437"
$ws.Range("E48").Value = "'-3.19%"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("E50").Value = "'-0.02%"
